# Daily attendance processing - 2025-10-29 03:03:33
# Rotate the "Recorded By" (column G) name/email list for each affected row:
# the first entry in the comma-separated list is moved to the end of the list.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(2,3,4,5,6,8,10,11,12,13,14,15,17,18,19,20,21,22,24,29,30,31,32,33,35,37,38,39,40,41,42,44,45,46,47,48,49,51,56,57,58,59,60,62,64,65,66,67,68,69,71,72,73,74,75,76,78,83,84,85,86,87,88,89,90,93,95,96,97,99,102,109,110,111,112,113,114,115,116,119,121,122,123,125,128,135,136,137,138,139,140,141,142,145,147,148,149,151,154)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    $parts = $current -split ", "
    if ($parts.Length -gt 1) {
        $rotated = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ", "
        $cell.Value = $rotated
    }
}
